# Commit: "Se cambia la palabra Alumnos por Estudiantes"
#
# Slide 1 ("TT IQQ 2018 Camara y Garay" title slide) had several shapes
# nudged slightly (the title placeholder gained an explicit <a:xfrm>, and
# three rectangles were repositioned a few EMU) together with the word
# "Alumnos" being renamed to "Estudiantes" in the authors rectangle.
#
# Point values below are chosen so that, after PowerPoint's internal
# Single-precision (EMU = floor(float32(pt) * 12700)) round trip, the
# saved OOXML offsets/extents land on the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "Título 1" (ctrTitle placeholder) -- gains an explicit xfrm
# matching its current effective (layout-inherited) position/size.
$title = $s.Shapes.Item(1)
$title.Left   = 174.2149658203125
$title.Top    = 186.56890869140625
$title.Width  = 611.56982421875
$title.Height = 129.6300811767578

# Shape 4: "Rectángulo 37" ("Memoria para optar al Título...") -- offset only.
$rect37 = $s.Shapes.Item(4)
$rect37.Left = 149.12417602539062
$rect37.Top  = 355.3514404296875

# Shape 5: "Rectángulo 39" (authors block) -- offset + word change.
$rect39 = $s.Shapes.Item(5)
$rect39.Left = 552.7424926757812
$rect39.Top  = 317.4283752441406
$rect39.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Estudiantes"

# Shape 6: "Rectángulo 42" ("Iquique - Chile" / "Diciembre 2018") -- offset only (x unchanged).
$rect42 = $s.Shapes.Item(6)
$rect42.Left = 225.4271697998047
$rect42.Top  = 467.20703125
